$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 15:10"

# Full refreshed country statistics table (row|country|CasosTotales|NuevosCasos|CasosActivos|Recuperados|CasosCriticos|MuertesHoy|Muertes)
# The source data is already sorted descending by "Casos totales" (column B), matching the
# canonical ordering used by the workbook.
$countryData = @(
"4|Estados Unidos|1746335|532|490256|1153963|0|9|102116",
"5|Brasil|414661|0|166647|222317|0|0|25697",
"6|Rusia|379051|8371|150993|223916|0|174|4142",
"7|España|283849|0|196958|59773|0|0|27118",
"8|Reino Unido|267240|0|0|0|0|0|37460",
"9|Italia|231139|0|147101|50966|0|0|33072",
"10|Francia|182913|0|66584|87733|0|0|28596",
"11|Alemania|181895|0|163200|10162|0|0|8533",
"12|Turquia|159797|0|122793|32573|0|0|4431",
"13|India|159054|968|67929|86584|0|7|4541",
"14|Iran|143849|2258|112988|23234|0|63|7627",
"15|Peru|135905|0|56169|75753|0|0|3983",
"16|Canada|87519|0|46164|34590|0|0|6765",
"17|China|82995|2|78288|73|0|0|4634",
"18|Chile|82289|0|33540|47908|0|0|841",
"19|Arabia Saudita|80185|1644|54553|25191|0|16|441",
"20|Mexico|78023|3463|54383|15043|0|463|8597",
"21|Pakistan|61227|2076|20231|39736|0|35|1260",
"22|Belgica|57849|257|15572|32889|0|24|9388",
"23|Catar|50914|1967|15399|35482|0|3|33",
"24|Paises Bajos|45950|182|0|0|0|32|5903",
"25|Banglades|40321|2029|8425|31337|0|15|559",
"26|Bielorrusia|39858|902|16660|22979|0|5|219",
"27|Ecuador|38103|0|18425|16403|0|0|3275",
"28|Suecia|35727|639|4971|26490|0|46|4266",
"29|Singapur|33249|373|17276|15950|0|0|23",
"30|Emiratos Arabes Unidos|32532|563|16685|15589|0|3|258",
"31|Portugal|31596|304|18637|11590|0|13|1369",
"32|Suiza|30796|20|28300|579|0|0|1917",
"33|Sudafrica|25937|0|13451|11934|0|0|552",
"34|Irlanda|24803|0|22089|1083|0|0|1631",
"35|Indonesia|24538|687|6240|16802|0|23|1496",
"36|Kuwait|24112|845|8698|15229|0|10|185",
"37|Colombia|24104|0|6111|17190|0|0|803",
"38|Polonia|22600|127|10560|11010|0|2|1030",
"39|Ucrania|22382|477|8439|13274|0|11|669",
"40|Egipto|19666|0|5205|13645|0|0|816",
"41|Rumania|18791|197|12629|4933|0|2|1229",
"42|Israel|16809|16|14602|1926|0|0|281",
"43|Japon|16651|0|13973|1820|0|0|858",
"44|Austria|16628|37|15286|674|0|23|668",
"45|Republica Dominicana|15723|0|8790|6459|0|0|474",
"46|Filipinas|15588|539|3598|11069|0|17|921",
"47|Argentina|13933|0|4349|9084|0|0|500",
"48|Afganistan|13036|580|1209|11592|0|8|235",
"49|Panama|11728|0|7379|4034|0|0|315",
"50|Dinamarca|11512|32|10180|764|0|3|568",
"51|Corea del Sur|11344|79|10340|735|0|0|269",
"52|Serbia|11300|25|6438|4621|0|1|241",
"53|Barein|9977|285|5295|4667|0|0|15",
"54|Kazajistan|9576|272|4893|4646|0|0|37",
"55|Chequia|9103|17|6377|2409|0|0|317",
"56|Oman|9009|636|2177|6792|0|1|40",
"57|Argelia|8857|0|5129|3105|0|0|623",
"58|Nigeria|8733|0|2501|5978|0|0|254",
"59|Noruega|8401|0|7727|438|0|1|236",
"60|Armenia|8216|442|3287|4816|0|15|113",
"61|Bolivia|7768|632|689|6799|0|6|280",
"62|Marruecos|7636|35|5109|2325|0|0|202",
"63|Malasia|7629|10|6169|1345|0|0|115",
"64|Moldavia|7537|0|4123|3138|0|2|276",
"65|Ghana|7303|0|2412|4857|0|0|34",
"66|Australia|7150|11|6579|468|0|0|103",
"67|Finlandia|6743|51|5500|930|0|0|313",
"68|Camerun|5436|0|1996|3265|0|0|175",
"69|Irak|5135|0|2904|2056|0|0|175",
"70|Honduras|4640|239|506|3940|0|6|194",
"71|Azerbaiyan|4568|0|2897|1617|0|0|54",
"72|Sudan|4346|200|749|3402|0|11|195",
"73|Guatemala|4145|191|493|3584|0|5|68",
"74|Luxemburgo|4001|0|3791|100|0|0|110",
"75|Hungria|3816|23|1996|1311|0|4|509",
"76|Tayikistan|3563|139|1674|1842|0|0|47",
"77|Uzbekistan|3437|68|2685|738|0|0|14",
"78|Senegal|3348|95|1686|1623|0|1|39",
"79|Guinea|3275|0|1673|1582|0|0|20",
"80|Tailandia|3065|11|2945|63|0|0|57",
"81|Grecia|2903|0|1374|1356|0|0|173",
"82|Republica de Yibuti|2697|0|1185|1494|0|0|18",
"83|Consejo Danes para los Refugiados|2660|114|381|2210|0|1|69",
"84|Costa de Marfil|2556|0|1302|1223|0|0|31",
"85|Bulgaria|2477|17|965|1378|0|1|134",
"86|Bosnia y Herzegovina|2462|27|1781|528|0|2|153",
"87|Gabon|2319|0|631|1674|0|0|14",
"88|Croacia|2245|1|2051|92|0|1|102",
"89|El Salvador|2194|85|1002|1153|0|2|39",
"90|Republica de Macedonia|2077|38|1486|470|0|2|121",
"91|Cuba|1974|0|1724|168|0|0|82",
"92|Estonia|1851|11|1574|211|0|0|66",
"93|Islandia|1805|0|1792|3|0|0|10",
"94|Somalia|1731|0|265|1399|0|0|67",
"95|Lituania|1656|9|1193|395|0|2|68",
"96|Mayotte|1645|0|1314|311|0|0|20",
"97|Kenia|1618|147|421|1139|0|3|58",
"98|Kirguistan|1594|74|1066|512|0|0|16",
"99|Eslovaquia|1520|5|1332|160|0|0|28",
"100|Nueva Zelanda|1504|0|1474|8|0|1|22",
"101|Sri Lanka|1486|17|745|731|0|0|10",
"102|Eslovenia|1473|2|1356|9|0|0|108",
"103|Maldivas|1457|0|197|1255|0|0|5",
"104|Haiti|1320|146|22|1264|0|1|34",
"105|Venezuela|1245|0|302|932|0|0|11",
"106|Guinea-Bisau|1195|0|42|1146|0|0|7",
"107|Libano|1168|7|699|443|0|0|26",
"108|Mali|1116|0|632|414|0|0|70",
"109|Albania|1076|26|823|220|0|0|33",
"110|Tunez|1068|17|938|82|0|0|48",
"111|Hong Kong|1067|0|1035|28|0|0|4",
"112|Letonia|1061|4|741|296|0|1|24",
"113|Zambia|1057|0|779|271|0|0|7",
"114|Guinea Ecuatorial|1043|0|165|866|0|0|12",
"115|Nepal|1042|156|187|850|0|1|5",
"116|Sudan del Sur|994|0|6|978|0|0|10",
"117|Costa Rica|984|0|639|335|0|0|10",
"118|Niger|955|3|796|95|0|1|64",
"119|Republica de Chipre|939|0|594|328|0|0|17",
"120|Paraguay|884|0|392|481|0|0|11",
"121|Burkina Faso|845|0|672|120|0|0|53",
"122|Etiopia|831|100|191|633|0|1|7",
"123|Sierra Leona|812|30|361|406|0|0|45",
"124|Uruguay|803|0|650|131|0|0|22",
"125|Principado de Andorra|763|0|676|36|0|0|51",
"126|Nicaragua|759|0|370|354|0|0|35",
"127|Georgia|738|3|573|153|0|0|12",
"128|Jordania|720|0|486|225|0|0|9",
"129|Republica del Chad|715|0|359|292|0|0|64",
"130|Crucero|712|0|651|48|0|0|13",
"131|Republica de Africa Central|702|0|23|678|0|0|1",
"132|San Marino|670|3|322|306|0|0|42",
"133|Madagascar|656|44|154|500|0|0|2",
"134|Malta|616|4|501|108|0|0|7",
"135|Congo|571|0|161|391|0|0|19",
"136|Jamaica|569|5|279|281|0|0|9",
"137|Tanzania|509|0|183|305|0|0|21",
"138|Reunion|460|0|411|48|0|0|1",
"139|Santo Tome y Principe|443|0|68|363|0|0|12",
"140|Taiwan|441|0|420|14|0|0|7",
"141|Estado de Palestina|435|1|365|67|0|0|3",
"142|Guayana Francesa|406|0|150|255|0|0|1",
"143|Togo|395|0|183|199|0|0|13",
"144|Cabo Verde|390|0|155|231|0|0|4",
"145|Ruanda|346|0|245|101|0|0|0",
"146|Isla de Man|336|0|306|6|0|0|24",
"147|Mauricio|334|0|322|2|0|0|10",
"148|Vietnam|327|0|278|49|0|0|0",
"149|Montenegro|324|0|315|0|0|0|9",
"150|Mauritania|292|0|15|261|0|0|16",
"151|Uganda|281|0|69|212|0|0|0",
"152|Suazilandia|272|0|168|102|0|0|2",
"153|Liberia|269|3|144|98|0|0|27",
"154|Yemen|256|0|10|193|0|0|53",
"155|Mozambique|227|0|71|155|0|0|1",
"156|Benin|210|0|134|73|0|0|3",
"157|Birmania|206|0|126|74|0|0|6",
"158|Martinica|197|0|91|92|0|0|14",
"159|Islas Feroe|187|0|187|0|0|0|0",
"160|Mongolia|161|13|43|118|0|0|0",
"161|Guadalupe|161|0|115|32|0|0|14",
"162|Gibraltar|158|1|147|11|0|0|0",
"163|Brunei|141|0|138|1|0|0|2",
"164|Islas Caimanes|140|0|67|72|0|0|1",
"165|Guyana|139|0|67|61|0|0|11",
"166|Bermudas|139|0|91|39|0|0|9",
"167|Zimbabue|132|0|25|103|0|0|4",
"168|Camboya|124|0|122|2|0|0|0",
"169|Siria|121|0|43|74|0|0|4",
"170|Trinidad yTobago|116|0|108|0|0|0|8",
"171|Malaui|101|0|37|60|0|0|4",
"172|Aruba|101|0|97|1|0|0|3",
"173|Bahamas|100|0|46|43|0|0|11",
"174|Libia|99|0|40|55|0|0|4",
"175|Monaco|98|0|90|4|0|0|4",
"176|Barbados|92|0|76|9|0|0|7",
"177|Comoras|87|0|24|61|0|0|2",
"178|Liechtenstein|82|0|55|26|0|0|1",
"179|San Martin (Parte Holandesa)|77|0|60|2|0|0|15",
"180|Angola|71|0|18|49|0|0|4",
"181|Polinesia Francesa|60|0|60|0|0|0|0",
"182|Macao|45|0|45|0|0|0|0",
"183|Burundi|42|0|20|21|0|0|1",
"184|San Martin (Parte Francesa)|40|0|33|4|0|0|3",
"185|Puerto Rico|39|0|1|36|0|0|2",
"186|Eritrea|39|0|39|0|0|0|0",
"187|Botsuana|35|0|20|14|0|0|1",
"188|Guam|32|0|0|31|0|0|1",
"189|Butan|31|3|6|25|0|0|0",
"190|San Vicente y las Granadinas|25|7|14|11|0|0|0",
"191|Gambia|25|0|19|5|0|0|1",
"192|Antigua y Barbuda|25|0|19|3|0|0|3",
"193|Timor Oriental|24|0|24|0|0|0|0",
"194|Granada|23|0|18|5|0|0|0",
"195|Namibia|22|0|14|8|0|0|0",
"196|Laos|19|0|16|3|0|0|0",
"197|Curazao|18|0|14|3|0|0|1",
"198|Fiyi|18|0|15|3|0|0|0",
"199|Santa Lucia|18|0|18|0|0|0|0",
"200|Belice|18|0|16|0|0|0|2",
"201|Nueva Caledonia|18|0|18|0|0|0|0",
"202|Islas Virgenes de los Estados Unidos|17|0|0|17|0|0|0",
"203|Dominica|16|0|16|0|0|0|0",
"204|San Cristobal y Nieves|15|0|15|0|0|0|0",
"205|Groenlandia|13|0|11|2|0|0|0",
"206|Islas Malvinas|13|0|13|0|0|0|0",
"207|Santa Sede|12|0|2|10|0|0|0",
"208|Surinam|12|0|9|2|0|0|1",
"209|Islas Turcas y Caicos|12|0|10|1|0|0|1",
"210|Montserrat|11|0|10|0|0|0|1",
"211|Seychelles|11|0|11|0|0|0|0",
"212|Sahara Occidental|9|0|6|2|0|0|1",
"213|Papua Nueva Guinea|8|0|8|0|0|0|0",
"214|Islas Virgenes Britanicas|8|0|7|0|0|0|1",
"215|San Bartolome|6|0|6|0|0|0|0",
"216|Bonaire, San Eustaquio y Saba|6|0|6|0|0|0|0",
"217|Anguila|3|0|3|0|0|0|0",
"218|Lesoto|2|0|0|2|0|0|0",
"219|San Pedro y Miquelon|1|0|1|0|0|0|0"
)

foreach ($line in $countryData) {
    $parts = $line.Split("|")
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = $parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [double]$parts[4]
    $ws.Cells.Item($r, 5).Value = [double]$parts[5]
    $ws.Cells.Item($r, 6).Value = [double]$parts[6]
    $ws.Cells.Item($r, 7).Value = [double]$parts[7]
    $ws.Cells.Item($r, 8).Value = [double]$parts[8]
}
